$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (D1:G1) - insert Corequisites, Concurrent, Recommended before
# the existing "Terms Typically Offered" which moves from D1 to G1.
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"
$ws.Range("G1").Value = "Terms Typically Offered"

# Row 2 (ESM 90) updates
$ws.Range("C2").Value = "Appropriate Math Placement Level."
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "ESM 105."
$ws.Range("F2").Value = "NA"
$ws.Range("G2").Value = "SU "

# Row 3 (ESM 105) updates
$ws.Range("C3").Value = "Appropriate Math Placement Level."
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "ESM 90."
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "SU "
